$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC!row20
$ws_ALC.Cells.Item(20, 8).Value = 1000
$ws_ALC.Cells.Item(20, 9).Value = 1000
$ws_ALC.Cells.Item(20, 11).Value = 1000
$ws_ALC.Cells.Item(20, 13).Value = -770

# ALC!row35
$ws_ALC.Cells.Item(35, 8).Value = 1000
$ws_ALC.Cells.Item(35, 9).Value = 1000
$ws_ALC.Cells.Item(35, 11).Value = 1000
$ws_ALC.Cells.Item(35, 13).Value = -621

# ALC!row43
$ws_ALC.Cells.Item(43, 8).Value = 7952694
$ws_ALC.Cells.Item(43, 9).Value = 100001
$ws_ALC.Cells.Item(43, 10).Value = 9261476
$ws_ALC.Cells.Item(43, 11).Value = 100001
$ws_ALC.Cells.Item(43, 12).Value = 9261476
$ws_ALC.Cells.Item(43, 13).Value = -99932
$ws_ALC.Cells.Item(43, 14).Value = -9261614

# ALC!row137
$ws_ALC.Cells.Item(137, 8).Value = 1552.9412
$ws_ALC.Cells.Item(137, 10).Value = 1789.091
$ws_ALC.Cells.Item(137, 12).Value = 5367.272999999999
$ws_ALC.Cells.Item(137, 14).Value = -10467.273

# ARM!row2
$ws_ARM.Cells.Item(2, 8).Value = 926.5714
$ws_ARM.Cells.Item(2, 9).Value = 757.4286
$ws_ARM.Cells.Item(2, 10).Value = 1095.7142
$ws_ARM.Cells.Item(2, 11).Value = 757.4286
$ws_ARM.Cells.Item(2, 12).Value = 1095.7142
$ws_ARM.Cells.Item(2, 13).Value = -644.4286
$ws_ARM.Cells.Item(2, 14).Value = -1321.7142

# ARM!row32
$ws_ARM.Cells.Item(32, 8).Value = 5724.6665
$ws_ARM.Cells.Item(32, 9).Value = 5724.6665
$ws_ARM.Cells.Item(32, 11).Value = 5724.6665
$ws_ARM.Cells.Item(32, 13).Value = -5437.6665

# ARM!row39
$ws_ARM.Cells.Item(39, 8).Value = 3900
$ws_ARM.Cells.Item(39, 9).Value = 3900
$ws_ARM.Cells.Item(39, 11).Value = 3900
$ws_ARM.Cells.Item(39, 13).Value = -3380

# ARM!row41
$ws_ARM.Cells.Item(41, 8).Value = 8292.666999999999
$ws_ARM.Cells.Item(41, 9).Value = 5689
$ws_ARM.Cells.Item(41, 10).Value = 13500
$ws_ARM.Cells.Item(41, 11).Value = 5689
$ws_ARM.Cells.Item(41, 12).Value = 13500
$ws_ARM.Cells.Item(41, 13).Value = -5275
$ws_ARM.Cells.Item(41, 14).Value = -14328

# ARM!row49
$ws_ARM.Cells.Item(49, 8).Value = 9500
$ws_ARM.Cells.Item(49, 10).Value = 9500
$ws_ARM.Cells.Item(49, 12).Value = 9500
$ws_ARM.Cells.Item(49, 14).Value = -10020

# ARM!row116
$ws_ARM.Cells.Item(116, 8).Value = 926.5714
$ws_ARM.Cells.Item(116, 9).Value = 757.4286
$ws_ARM.Cells.Item(116, 10).Value = 1095.7142
$ws_ARM.Cells.Item(116, 11).Value = 757.4286
$ws_ARM.Cells.Item(116, 12).Value = 1095.7142
$ws_ARM.Cells.Item(116, 13).Value = 1536.5714
$ws_ARM.Cells.Item(116, 14).Value = -5683.7142

# ARM!row132
$ws_ARM.Cells.Item(132, 8).Value = 3129.7693
$ws_ARM.Cells.Item(132, 9).Value = 3002.25
$ws_ARM.Cells.Item(132, 10).Value = 3333.8
$ws_ARM.Cells.Item(132, 11).Value = 9006.75
$ws_ARM.Cells.Item(132, 12).Value = 10001.4
$ws_ARM.Cells.Item(132, 13).Value = -6476.75
$ws_ARM.Cells.Item(132, 14).Value = -15061.4

# BSM!row3
$ws_BSM.Cells.Item(3, 8).Value = 926.5714
$ws_BSM.Cells.Item(3, 9).Value = 757.4286
$ws_BSM.Cells.Item(3, 10).Value = 1095.7142
$ws_BSM.Cells.Item(3, 11).Value = 757.4286
$ws_BSM.Cells.Item(3, 12).Value = 1095.7142
$ws_BSM.Cells.Item(3, 13).Value = -643.4286
$ws_BSM.Cells.Item(3, 14).Value = -1323.7142

# BSM!row38
$ws_BSM.Cells.Item(38, 8).Value = 0
$ws_BSM.Cells.Item(38, 10).Value = 0
$ws_BSM.Cells.Item(38, 12).Value = 0
$ws_BSM.Cells.Item(38, 14).ClearContents()

# BSM!row56
$ws_BSM.Cells.Item(56, 8).Value = 0
$ws_BSM.Cells.Item(56, 10).Value = 0
$ws_BSM.Cells.Item(56, 12).Value = 0
$ws_BSM.Cells.Item(56, 14).ClearContents()

# BSM!row94
$ws_BSM.Cells.Item(94, 8).Value = 62500224
$ws_BSM.Cells.Item(94, 9).Value = 62500224
$ws_BSM.Cells.Item(94, 11).Value = 62500224
$ws_BSM.Cells.Item(94, 13).Value = -62499773

# BSM!row99
$ws_BSM.Cells.Item(99, 8).Value = 125001160
$ws_BSM.Cells.Item(99, 9).Value = 200001090
$ws_BSM.Cells.Item(99, 11).Value = 200001090
$ws_BSM.Cells.Item(99, 13).Value = -199999592

# CRP!row16
$ws_CRP.Cells.Item(16, 8).Value = 333334100
$ws_CRP.Cells.Item(16, 9).Value = 333334100
$ws_CRP.Cells.Item(16, 10).Value = 0
$ws_CRP.Cells.Item(16, 11).Value = 333334100
$ws_CRP.Cells.Item(16, 12).Value = 0
$ws_CRP.Cells.Item(16, 13).Value = -333333813
$ws_CRP.Cells.Item(16, 14).ClearContents()

# CRP!row31
$ws_CRP.Cells.Item(31, 8).Value = 1239.6383
$ws_CRP.Cells.Item(31, 9).Value = 1271.7142
$ws_CRP.Cells.Item(31, 10).Value = 1226.0303
$ws_CRP.Cells.Item(31, 11).Value = 1271.7142
$ws_CRP.Cells.Item(31, 12).Value = 1226.0303
$ws_CRP.Cells.Item(31, 13).Value = -976.7141999999999
$ws_CRP.Cells.Item(31, 14).Value = -1816.0303

# CRP!row34
$ws_CRP.Cells.Item(34, 8).Value = 1239.6383
$ws_CRP.Cells.Item(34, 9).Value = 1271.7142
$ws_CRP.Cells.Item(34, 10).Value = 1226.0303
$ws_CRP.Cells.Item(34, 11).Value = 1271.7142
$ws_CRP.Cells.Item(34, 12).Value = 1226.0303
$ws_CRP.Cells.Item(34, 13).Value = -1069.7142
$ws_CRP.Cells.Item(34, 14).Value = -1630.0303

# CRP!row58
$ws_CRP.Cells.Item(58, 8).Value = 1571
$ws_CRP.Cells.Item(58, 9).Value = 1224
$ws_CRP.Cells.Item(58, 10).Value = 2314.5715
$ws_CRP.Cells.Item(58, 11).Value = 1224
$ws_CRP.Cells.Item(58, 12).Value = 2314.5715
$ws_CRP.Cells.Item(58, 13).Value = -1021
$ws_CRP.Cells.Item(58, 14).Value = -2720.5715

# CRP!row69
$ws_CRP.Cells.Item(69, 8).Value = 0
$ws_CRP.Cells.Item(69, 9).Value = 0
$ws_CRP.Cells.Item(69, 11).Value = 0
$ws_CRP.Cells.Item(69, 13).ClearContents()

# CRP!row72
$ws_CRP.Cells.Item(72, 8).Value = 0
$ws_CRP.Cells.Item(72, 9).Value = 0
$ws_CRP.Cells.Item(72, 11).Value = 0
$ws_CRP.Cells.Item(72, 13).ClearContents()

# CRP!row99
$ws_CRP.Cells.Item(99, 8).Value = 1528.3636
$ws_CRP.Cells.Item(99, 9).Value = 1464.2858
$ws_CRP.Cells.Item(99, 10).Value = 1640.5
$ws_CRP.Cells.Item(99, 11).Value = 1464.2858
$ws_CRP.Cells.Item(99, 12).Value = 1640.5
$ws_CRP.Cells.Item(99, 13).Value = 33.71419999999989
$ws_CRP.Cells.Item(99, 14).Value = -4636.5

# CRP!row105
$ws_CRP.Cells.Item(105, 8).Value = 700
$ws_CRP.Cells.Item(105, 9).Value = 700
$ws_CRP.Cells.Item(105, 10).Value = 0
$ws_CRP.Cells.Item(105, 11).Value = 700
$ws_CRP.Cells.Item(105, 12).Value = 0
$ws_CRP.Cells.Item(105, 13).Value = 1047
$ws_CRP.Cells.Item(105, 14).ClearContents()

# CRP!row107
$ws_CRP.Cells.Item(107, 8).Value = 706.44446
$ws_CRP.Cells.Item(107, 9).Value = 339
$ws_CRP.Cells.Item(107, 10).Value = 1992.5
$ws_CRP.Cells.Item(107, 11).Value = 339
$ws_CRP.Cells.Item(107, 12).Value = 1992.5
$ws_CRP.Cells.Item(107, 13).Value = 1581
$ws_CRP.Cells.Item(107, 14).Value = -5832.5

# CRP!row113
$ws_CRP.Cells.Item(113, 8).Value = 333334100
$ws_CRP.Cells.Item(113, 9).Value = 333334100
$ws_CRP.Cells.Item(113, 10).Value = 0
$ws_CRP.Cells.Item(113, 11).Value = 333334100
$ws_CRP.Cells.Item(113, 12).Value = 0
$ws_CRP.Cells.Item(113, 13).Value = -333331930
$ws_CRP.Cells.Item(113, 14).ClearContents()

# CRP!row122
$ws_CRP.Cells.Item(122, 8).Value = 650.2
$ws_CRP.Cells.Item(122, 9).Value = 666.2941
$ws_CRP.Cells.Item(122, 10).Value = 559
$ws_CRP.Cells.Item(122, 11).Value = 1998.8823
$ws_CRP.Cells.Item(122, 12).Value = 1677
$ws_CRP.Cells.Item(122, 13).Value = 451.1177000000002
$ws_CRP.Cells.Item(122, 14).Value = -6577

# CRP!row126
$ws_CRP.Cells.Item(126, 8).Value = 1528.3636
$ws_CRP.Cells.Item(126, 9).Value = 1464.2858
$ws_CRP.Cells.Item(126, 10).Value = 1640.5
$ws_CRP.Cells.Item(126, 11).Value = 4392.857400000001
$ws_CRP.Cells.Item(126, 12).Value = 4921.5
$ws_CRP.Cells.Item(126, 13).Value = -1922.857400000001
$ws_CRP.Cells.Item(126, 14).Value = -9861.5

# CRP!row134
$ws_CRP.Cells.Item(134, 8).Value = 15153118
$ws_CRP.Cells.Item(134, 9).Value = 1681.0385
$ws_CRP.Cells.Item(134, 10).Value = 71429880
$ws_CRP.Cells.Item(134, 11).Value = 5043.1155
$ws_CRP.Cells.Item(134, 12).Value = 214289640
$ws_CRP.Cells.Item(134, 13).Value = -2508.1155
$ws_CRP.Cells.Item(134, 14).Value = -214294710

# CRP!row136
$ws_CRP.Cells.Item(136, 8).Value = 1571
$ws_CRP.Cells.Item(136, 9).Value = 1224
$ws_CRP.Cells.Item(136, 10).Value = 2314.5715
$ws_CRP.Cells.Item(136, 11).Value = 3672
$ws_CRP.Cells.Item(136, 12).Value = 6943.7145
$ws_CRP.Cells.Item(136, 13).Value = -1122
$ws_CRP.Cells.Item(136, 14).Value = -12043.7145

# CUL!row15
$ws_CUL.Cells.Item(15, 8).Value = 322.22223
$ws_CUL.Cells.Item(15, 9).Value = 322.22223
$ws_CUL.Cells.Item(15, 11).Value = 966.66669
$ws_CUL.Cells.Item(15, 13).Value = -826.66669

# CUL!row20
$ws_CUL.Cells.Item(20, 8).Value = 300
$ws_CUL.Cells.Item(20, 9).Value = 300
$ws_CUL.Cells.Item(20, 10).Value = 300
$ws_CUL.Cells.Item(20, 11).Value = 900
$ws_CUL.Cells.Item(20, 12).Value = 900
$ws_CUL.Cells.Item(20, 13).Value = -673
$ws_CUL.Cells.Item(20, 14).Value = -1354

# GSM!row2
$ws_GSM.Cells.Item(2, 8).Value = 225.1
$ws_GSM.Cells.Item(2, 9).Value = 170.8
$ws_GSM.Cells.Item(2, 11).Value = 170.8
$ws_GSM.Cells.Item(2, 13).Value = -57.80000000000001

# GSM!row102
$ws_GSM.Cells.Item(102, 8).Value = 1274.8572
$ws_GSM.Cells.Item(102, 9).Value = 1248.2222
$ws_GSM.Cells.Item(102, 10).Value = 1322.8
$ws_GSM.Cells.Item(102, 11).Value = 1248.2222
$ws_GSM.Cells.Item(102, 12).Value = 1322.8
$ws_GSM.Cells.Item(102, 13).Value = 373.7778000000001
$ws_GSM.Cells.Item(102, 14).Value = -4566.8

# GSM!row113
$ws_GSM.Cells.Item(113, 8).Value = 1282.8572
$ws_GSM.Cells.Item(113, 9).Value = 1140
$ws_GSM.Cells.Item(113, 10).Value = 1540
$ws_GSM.Cells.Item(113, 11).Value = 1140
$ws_GSM.Cells.Item(113, 12).Value = 1540
$ws_GSM.Cells.Item(113, 13).Value = 1030
$ws_GSM.Cells.Item(113, 14).Value = -5880

# GSM!row122
$ws_GSM.Cells.Item(122, 8).Value = 7814128.5
$ws_GSM.Cells.Item(122, 9).Value = 1632.5238
$ws_GSM.Cells.Item(122, 11).Value = 4897.5714
$ws_GSM.Cells.Item(122, 13).Value = -2447.5714

# GSM!row132
$ws_GSM.Cells.Item(132, 8).Value = 3229.3333
$ws_GSM.Cells.Item(132, 9).Value = 4009
$ws_GSM.Cells.Item(132, 10).Value = 2520.5454
$ws_GSM.Cells.Item(132, 11).Value = 12027
$ws_GSM.Cells.Item(132, 12).Value = 7561.6362
$ws_GSM.Cells.Item(132, 13).Value = -9497
$ws_GSM.Cells.Item(132, 14).Value = -12621.6362

# LTW!row22
$ws_LTW.Cells.Item(22, 8).Value = 693.0714
$ws_LTW.Cells.Item(22, 9).Value = 465.2857
$ws_LTW.Cells.Item(22, 11).Value = 465.2857
$ws_LTW.Cells.Item(22, 13).Value = -170.2857

# LTW!row27
$ws_LTW.Cells.Item(27, 8).Value = 693.0714
$ws_LTW.Cells.Item(27, 9).Value = 465.2857
$ws_LTW.Cells.Item(27, 11).Value = 465.2857
$ws_LTW.Cells.Item(27, 13).Value = -358.2857

# LTW!row68
$ws_LTW.Cells.Item(68, 8).Value = 1705.5238
$ws_LTW.Cells.Item(68, 9).Value = 1712.0555
$ws_LTW.Cells.Item(68, 11).Value = 1712.0555
$ws_LTW.Cells.Item(68, 13).Value = -963.0554999999999

# LTW!row71
$ws_LTW.Cells.Item(71, 8).Value = 1705.5238
$ws_LTW.Cells.Item(71, 9).Value = 1712.0555
$ws_LTW.Cells.Item(71, 11).Value = 8560.2775
$ws_LTW.Cells.Item(71, 13).Value = -4816.2775

# LTW!row93
$ws_LTW.Cells.Item(93, 8).Value = 1015.4
$ws_LTW.Cells.Item(93, 9).Value = 991.6667
$ws_LTW.Cells.Item(93, 10).Value = 1051
$ws_LTW.Cells.Item(93, 11).Value = 991.6667
$ws_LTW.Cells.Item(93, 12).Value = 1051
$ws_LTW.Cells.Item(93, 13).Value = 256.3333
$ws_LTW.Cells.Item(93, 14).Value = -3547

# LTW!row100
$ws_LTW.Cells.Item(100, 8).Value = 1320.5714
$ws_LTW.Cells.Item(100, 10).Value = 1338.75
$ws_LTW.Cells.Item(100, 12).Value = 1338.75
$ws_LTW.Cells.Item(100, 14).Value = -2420.75

# LTW!row122
$ws_LTW.Cells.Item(122, 8).Value = 19232646
$ws_LTW.Cells.Item(122, 9).Value = 35715970
$ws_LTW.Cells.Item(122, 10).Value = 2100.6667
$ws_LTW.Cells.Item(122, 11).Value = 107147910
$ws_LTW.Cells.Item(122, 12).Value = 6302.000100000001
$ws_LTW.Cells.Item(122, 13).Value = -107145460
$ws_LTW.Cells.Item(122, 14).Value = -11202.0001

# WVR!row122
$ws_WVR.Cells.Item(122, 8).Value = 12500907
$ws_WVR.Cells.Item(122, 9).Value = 14706832
$ws_WVR.Cells.Item(122, 10).Value = 665
$ws_WVR.Cells.Item(122, 11).Value = 44120496
$ws_WVR.Cells.Item(122, 12).Value = 1995
$ws_WVR.Cells.Item(122, 13).Value = -44118046
$ws_WVR.Cells.Item(122, 14).Value = -6895

# WVR!row126
$ws_WVR.Cells.Item(126, 8).Value = 125001070
$ws_WVR.Cells.Item(126, 9).Value = 142858080
$ws_WVR.Cells.Item(126, 11).Value = 428574240
$ws_WVR.Cells.Item(126, 13).Value = -428571770

# WVR!row136
$ws_WVR.Cells.Item(136, 8).Value = 979.1429000000001
$ws_WVR.Cells.Item(136, 9).Value = 919.08
$ws_WVR.Cells.Item(136, 11).Value = 2757.24
$ws_WVR.Cells.Item(136, 13).Value = -207.2400000000002
